$d = $word.ActiveDocument

# 1) Tighten the summary wording: "to the roll out of quality" -> "to roll out quality"
$ok1 = $d.Content.Find.Execute(
    "to the roll out of quality", $true, $false, $false, $false, $false,
    $true, 1, $false, "to roll out quality", 2)

# 2) Swap the closing phrase: "to broaden all skills." -> "to overcome any challenge."
$ok2 = $d.Content.Find.Execute(
    "to broaden all skills.", $true, $false, $false, $false, $false,
    $true, 1, $false, "to overcome any challenge.", 2)

# 3) Move the (hidden) "_GoBack" bookmark from the end of the "...C#" bullet
#    paragraph up to the end of the rewritten summary paragraph. Since
#    Word only ever keeps a single "_GoBack" bookmark, re-adding it with
#    that name anywhere else automatically removes the old one.
#
#    A truly collapsed range sitting right at the end of a paragraph (just
#    before the paragraph mark) gets snapped to the start of that
#    paragraph by this host's Bookmarks.Add, so we briefly insert a
#    placeholder character, bookmark around it, then delete the
#    placeholder - leaving a zero-length bookmark in the right spot, i.e.
#    exactly after the last run of text and before the paragraph mark.
$found = $d.Content
$found.Find.Execute("to overcome any challenge.", $true) | Out-Null
$target = $found.Duplicate
$target.Start = $target.End
$target.InsertAfter("x")
$target.Bookmarks.Add("_GoBack")
$target.Text = ""

Write-Output "replace1=$ok1 replace2=$ok2"
